$wb = $excel.ActiveWorkbook

# --- Sheet1: insert a new "ExpectedStatuscode" column between the existing
#     Input/Operation-Type columns and the ResponseParam/validation columns ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new column at G; existing G (ResponseParam) and H (validation) shift to H and I
$ws1.Columns("G").Insert()

# Populate the new header and value
$ws1.Range("G1").Value = "ExpectedStatuscode"
$ws1.Range("G2").Value = 200

# Match the width of the neighbouring column for the newly inserted column
$ws1.Columns("G").ColumnWidth = $ws1.Columns("F").ColumnWidth

# --- Sheet2 ("FP Flow"): just a change of the active selection ---
$ws2 = $wb.Worksheets.Item("FP Flow")
$ws2.Activate() | Out-Null
$ws2.Range("B4").Select() | Out-Null

# Re-activate Sheet1 (it remains the tab that is selected/visible) and
# update its active selection to the new status-code cell
$ws1.Activate() | Out-Null
$ws1.Range("G2").Select() | Out-Null
